$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.22222222222222
$ws.Range("C2").Value = 0.3703703703703703
$ws.Range("D2").Value = 99.66666666666667
$ws.Range("E2").Value = 30
$ws.Range("F2").Value = 5.259259259259259
$ws.Range("G2").Value = 34.96296296296296
$ws.Range("H2").Value = 69.14814814814814
$ws.Range("I2").Value = 19.92592592592593
$ws.Range("J2").Value = 21.07407407407408
$ws.Range("K2").Value = 26.13121201216228
